$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.286.47"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.16%  "

$ws.Range("D3").Value = "'1.858.80"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.27%  "

$ws.Range("D4").Value = "'1.017"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +1.19%  "

$ws.Range("D5").Value = "'314.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.26%  "

$ws.Range("D6").Value = "'1.013"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.94%  "

$ws.Range("D7").Value = "'0.5105"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.64%  "

$ws.Range("D8").Value = "'0.3922"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.45%  "

$ws.Range("D9").Value = "'0.08277"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.22%  "

$ws.Range("D10").Value = "'1.112"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.69%  "

$ws.Range("D11").Value = "'6.210"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.57%  "

$ws.Range("D12").Value = "'1.878.02"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.43%  "

$ws.Range("D13").Value = "'20.25"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.51%  "

$ws.Range("D14").Value = "'7.192"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.53%  "

$ws.Range("D15").Value = "'1.016"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.18%  "

$ws.Range("D16").Value = "'0.00001099"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.84%  "

$ws.Range("D17").Value = "'91.04"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.49%  "

$ws.Range("D18").Value = "'0.06692"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.42%  "

$ws.Range("D19").Value = "'17.57"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.25%  "

$ws.Range("D20").Value = "'1.014"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.99%  "

$ws.Range("D21").Value = "'5.935"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.11%  "

$ws.Range("D22").Value = "'28.319.94"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.17%  "

$ws.Range("D23").Value = "'11.04"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.56%  "

$ws.Range("D24").Value = "'2.257"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.52%  "

$ws.Range("D25").Value = "'2.071.77"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.00%  "

$ws.Range("D26").Value = "'160.35"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.86%  "

$ws.Range("D27").Value = "'20.56"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.53%  "

$ws.Range("D28").Value = "'2.389"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.05%  "

$ws.Range("D29").Value = "'126.42"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.64%  "

$ws.Range("D30").Value = "'0.1048"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.76%  "

$ws.Range("E31").Value = "  -1.12%  "

$ws.Range("D32").Value = "'5.788"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.89%  "

$ws.Range("D33").Value = "'3.629"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.97%  "

$ws.Range("D34").Value = "'0.02426"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.57%  "

$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "'0.06454"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.66%  "

$ws.Range("B36").Value = "FraxShare"
$ws.Range("C36").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D36").Value = "'9.128"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.49%  "

$ws.Range("D37").Value = "'0.2165"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.42%  "

$ws.Range("D38").Value = "'1.252"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.67%  "

$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").Value = "'1.180"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.63%  "

$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "'0.6406"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.99%  "

$ws.Range("D41").Value = "'4.936"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.06%  "

$ws.Range("D42").Value = "'11.09"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.19%  "

$ws.Range("D43").Value = "'0.5981"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.42%  "

$ws.Range("E44").Value = "  -2.10%  "

$ws.Range("D45").Value = "'3.689"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.20%  "

$ws.Range("D46").Value = "'1.282"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.81%  "

$ws.Range("D47").Value = "'1.976"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.39%  "

$ws.Range("D48").Value = "'1.199"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.77%  "

$ws.Range("D49").Value = "'120.44"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.15%  "

$ws.Range("D50").Value = "'0.06861"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.71%  "

$ws.Range("E51").Value = "  -3.44%  "
